# Apply "hybrid bold + color" quantitative-metrics highlighting to six
# bullet/impact lines in the resume. For each target paragraph we locate
# the full paragraph text with Find.Execute (searching forward from a
# cursor so repeated/prefix text is matched in document order), then for
# every metric substring (in left-to-right order) we compute its offset
# inside the matched text with .IndexOf and bold+color that sub-Range.
# Word's object model automatically re-splits the underlying w:r runs at
# the Range boundaries, which reproduces the run structure in the diff.

$d = $word.ActiveDocument

# Word's Font.Color is a BGR-packed long (not RGB), so 0x2C3E50 (the
# diff's <w:color w:val="2C3E50"/>) must be byte-swapped to 0x503E2C.
$metricColor = 0x503E2C

# Search cursor: advances past each match so a second occurrence of the
# same (or a prefix) string is found after the first, not re-matched.
$script:searchCursor = 0

function Apply-MetricHighlight($paragraphText, $metrics) {
    $docEnd = $d.Content.End
    $searchRange = $d.Range($script:searchCursor, $docEnd)
    $found = $searchRange.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Paragraph text not found: $paragraphText"
    }
    $pStart = $searchRange.Start
    $pEnd = $searchRange.End
    $script:searchCursor = $pEnd

    $searchFrom = 0
    foreach ($metric in $metrics) {
        $idx = $paragraphText.IndexOf($metric, $searchFrom)
        if ($idx -lt 0) {
            throw "Metric '$metric' not found in paragraph text"
        }
        $metricStart = $pStart + $idx
        $metricEnd = $metricStart + $metric.Length
        $metricRange = $d.Range($metricStart, $metricEnd)
        $metricRange.Font.Bold = 1
        $metricRange.Font.Color = $metricColor
        $searchFrom = $idx + $metric.Length
    }
}

# 1. Race-coding-errors achievement bullet (Siege Analytics)
Apply-MetricHighlight `
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" `
    @("23%", "64%")

# 2. Prediction-accuracy bullet with polling error margins (Siege Analytics)
Apply-MetricHighlight `
    "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%" `
    @("87%", "71%", "±4.2%", "±2.1%")

# 3. RFP vendor bids bullet (Myers Research)
Apply-MetricHighlight `
    "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development" `
    @("1,200")

# 4. Polling Consortium Database bullet (Lake Research Partners)
Apply-MetricHighlight `
    "• Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+" `
    @("`$400M", "`$1B")

# 5. Mapping-cost-reduction bullet (Key Achievements and Impact)
Apply-MetricHighlight `
    "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M" `
    @("73.5%", "`$4.7M")

# 6. Prediction-accuracy bullet without polling margins (Key Achievements and Impact)
Apply-MetricHighlight `
    "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" `
    @("87%", "71%")
